$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G ("d=7" column), shifting the
# existing "d=7" and "d=10" columns (and their data) one column to the right.
$ws.Columns("G:G").Insert()

# New header for the inserted column (copy the header formatting from
# the neighboring "d=5" header cell so the new cell gets the same style)
$ws.Range("G1").Value = "d=6"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data values for the inserted column
$ws.Range("G2").Value = 97.71483630067273
$ws.Range("G3").Value = 97.8415916170577
$ws.Range("G4").Value = 97.77716205695562
$ws.Range("G5").Value = 97.71973125282346
$ws.Range("G6").Value = 97.76820069040014
